$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NBS_CESY_process")

# Insert two new rows right after the title row (row 1), pushing all
# existing data down by two rows.
$ws.Rows.Item(2).Resize(2).Insert()

$ws.Range("A3").Value = "# ----------"
$ws.Range("A2").Value = "# Column types: cc"

$ws.Range("A2").Select()
